$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) D29 ("서비스 컨택 포인트" detail, HT/CP00002): split the "/"-joined
#    branch contact lines into individual "<br>" lines, and fix the
#    "jongwoon.yun@volvoty.co.krr" typo -> "...co.kr"
# ---------------------------------------------------------------------------
$d29 = @'
<a style="color: red; font-weight: bold;">Volvo Internal Use Only</a>
<br> 서비스센터를 통해 문의가 필요한 경우 아래 주소를 참고합니다. 유선 컨택 포인트는 액셀 파일을 참고하십시오.
<br>
<br><h3>에이치모터스</h3>
<br><string>본부장: 신홍열</string> hongyeul.shin@hvolvo.com
<br><string>CR매니저/본부장 서포트: 김진용</string> jinyong.kim@hvolvo.com
<br><string>서비스센터 (지점장):</string>
<br>강남 삼성 박성우: sungwoo.park@hvolvo.com
<br>성수 주상표: sangpyo.joo@hvolvo.com
<br>분당 서현 안정식: ahn.jeongsik@hvolvo.com
<br>율현 황의섭: euiseop.hwang@hvolvo.com
<br>수원 문정규: jeongkyu.moon@hvolvo.com
<br>인천 안종진: chongjin.an@hvolvo.com
<br>대전 전용일: yongil.jeon@hvolvo.com
<br>청주 고재국: jaeguk.ko@hvolvo.com
<br>
<br><h3>아이비모터스</h3>
<br><string>본부장: 임용혁<string> yonghyuk.im@ivymotors.co.kr
<br><string>서비스센터 (지점장):</string>
<br>광주 조요근: yokeun.jo@ivymotors.co.kr
<br>전주 황금용: geumyoug.hwang@ivymotors.co.kr
<br>순천 문경철: kyungchul.mun@ivymotors.co.kr
<br>제주 안대환: daehwan.an@ivymotors.co.kr
<br>
<br><h3>태영모터스</h3>
<br><string>본부장: 윤종운</string> jongwoon.yun@volvoty.co.kr
<br><string>CR매니저/본부장 서포트: 윤우식</string> woosik.yoon@volvoty.co.kr 
<br><string>서비스센터 (지점장):</string>
<br>대구 윤종운 본부장: jongwoon.yun@volvoty.co.kr
<br>포항 윤종운 본부장: jongwoon.yun@volvoty.co.kr
<br>서대구 권회열: hoiyol.kwon@volvoty.co.kr
<br>
<br><h3>천하자동차</h3>
<br><string>본부장: 최병재</string> byungjae.choi@chvolvo.co.kr
<br><string>서비스센터 (지점장):</string>
<br>영등포 양평 이학제: hakje.lee@chvolvo.co.kr
<br>동대문 김병창: byungchang.kim@chvolvo.co.kr
<br>의정부 김혁진 팀장: hyuckjin.kim@chvolvo.co.kr
<br>구리 이종모 선임: jongmo.lee@chvolvo.co.kr
<br>
<br><h3>아이언모터스</h3>
<br><string>본부장: 김인호</string> Inho.kim@ironmotors.co.kr
<br><string>서비스센터 (지점장):</string>
<br>창원 허명욱: myunguk.heo@ironmotors.co.kr
<br>김해 장성용: sungyong.jang@ironmotors.co.kr
<br>광안 조재우: jaewoo.cho@ironmotors.co.kr
<br>해운대 김영곤: younggon.kim@ironmotors.co.kr
<br>
<br><h3>아주오토리움</h3>
<br><string>본부장: 서일선</string> isseo@aju.co.kr
<br><string>서비스센터 (지점장):</string>
<br>일산 홍영기: yghong@aju.co.kr
<br>안양 이조원: jwlee8@aju.co.kr
<br>영등포 문래 서일선: isseo@aju.co.kr
<br>부천 이재석: leejs@aju.co.kr
<br>
<br><h3>코오롱 오토모티브</h3>
<br><string>본부장: 김태길</string> taegil_kim@kolon.com
<br><string>서비스센터 (지점장):</string>
<br>송파 정전성: jeonseong_jeong@kolon.com
<br>서초 최경모: kyoungmo_choi@kolon.com
<br>원주 심재섭: jaeseop_shim@kolon.com
<br>천안 최락천: rockchun_choi@kolon.com
<br>울산 임기범: Kibum_lim@kolon.com
<br>분당 판교 김종호: jongho_kim@kolon.com
<br>강릉 김태길: taegil_kim@kolon.com
<br>서산 최락천: rockchun_choi@kolon.com
<br>하남 김성준: sungjoon_kim@kolon.com
'@
$ws.Range("D29").Value = $d29

# ---------------------------------------------------------------------------
# 2) D30 ("세일즈 컨택 포인트" detail, CP00003): same "/"-split treatment,
#    plus the opening <a style="..."> tag gets its quotes doubled
#    ( style="..."  ->  style=""..."" ) exactly as in the source edit.
# ---------------------------------------------------------------------------
$d30 = @'
<a style=""color: red; font-weight: bold;"">Volvo Internal Use Only</a>
<br> 전시장을 통해 문의가 필요한 경우 아래 주소를 참고합니다. 유선 컨택 포인트는 액셀 파일을 참고하십시오.
<br>
<br><h3>에이치모터스</h3>
<br><string>대표: 황호진</string> ceo@hvolvo.com
<br><string>본부장: 예정규</string> jkye@hvolvo.com
<br><string>전시장 (지점장):</string>
<br>강남대치: 김길성 kskim@hvolvo.com
<br>강남신사: 강은식 eunsik.kang@hvolvo.com
<br>분당서현 양승혁:  shyang@hvolvo.com
<br>수원 유상곤: sanggon.you@hvolvo.com
<br>인천 홍동현: dhhong@hvolvo.com
<br>대전 이상호: leesh@hvolvo.com
<br>청주 윤상무: sangmoo.yun@hvolvo.com
<br>SELEKT 수원 박시현: sihyun.park@hvolvo.com
<br>
<br><h3>아이비모터스</h3>
<br><string>대표: 강병철</string> byungcheol.kang@ivymotors.co.kr
<br><string>본부장: 배성일<string> sungil.bae@ivymotors.co.kr
<br><string>전시장 (지점장):</string>
<br>광주 손영: young.son@ivymotors.co.kr
<br>전주 임동진: hwanwook.shin@ivymotors.co.kr
<br>순천 임동진: dongjin.lim@ivymotors.co.kr
<br>제주 전명호: myeongho.cheon@ivymotors.co.kr
<br>
<br><h3>태영모터스</h3>
<br><string>대표: 김용수</string> ys_kim@volvoty.co.kr
<br><string>본부장: 김상원</string> swkim@volvoty.co.kr
<br><string>전시장 (지점장):</string>
<br>대구 신승욱: sw.shin72@volvoty.co.kr
<br>포항 이경식: kslee_ty@volvoty.co.kr
<br>서대구 황영상: votycar@volvoty.co.kr
<br>
<br><h3>천하자동차</h3>
<br><string>대표: 윤인경</string> jamescap2000@chvolvo.co.kr
<br><string>본부장: 장준원</string> jwc@chvolvo.co.kr
<br><string>전시장 (지점장):</string>
<br>동대문 허민철: mcheo@chvolvo.co.kr
<br>의정부 김석영: ysk@chvolvo.co.kr
<br>구리 임승현: shlim@chvolvo.co.kr
<br>용산 강명윤: myungyun.kang@chvolvo.co.kr
<br>
<br><h3>아이언모터스</h3>
<br><string>대표: 김민규</string> minkyu.kim@ironmotors.co.kr
<br><string>본부장: 신홍섭</string> hongsub.shin@ironmotors.co.kr
<br><string>전시장 (지점장):</string>
<br>창원 김희종: heejong.kim@ironmotors.co.kr
<br>광안 안위성: wiseong.ahn@ironmotors.co.kr
<br>해운대 박성준: sungjun.park@ironmotors.co.kr
<br>김해 최경철: gyeongchoel.choe@ironmotors.co.kr
<br>진주 안정수: jeongsoo.ahn@ironmotors.co.kr
<br>SELEKT 부산 김호영: hoyoung.kim@ironmotors.co.kr
<br>
<br><h3>아주오토리움</h3>
<br><string>대표: 박영석</string> andypark@aju.co.kr
<br><string>본부장: 한영수</string> youngsoohan@aju.co.kr
<br><string>전시장 (지점장):</string>
<br>목동 이상일: silee@aju.co.kr
<br>일산 한택주: taek6664@aju.co.kr
<br>안양 문상호: shmun@aju.co.kr
<br>부천 박용호: yhpark05@aju.co.kr
<br>고양 임지운: jwlim@aju.co.kr
<br>
<br><h3>코오롱 오토모티브</h3>
<br><string>대표: 신진욱</string> jinuk_shin@kolon.com
<br><string>본부장: 최형준</string> hyungjun_choi@kolon.com
<br><string>전시장 (지점장):</string>
<br>송파 유지훈: jihoon_yoo2@kolon.com
<br>서초 홍영삼: cofficer@kolon.com
<br>원주 박월준: woljun_park@kolon.com
<br>천안 김영민: youngmin_kim1@kolon.com
<br>울산 소상만: mynia@kolon.com
<br>분당판교 김영선: youngsun_kim1@kolon.com
<br>하남 김한묵: hanmook_kim@kolon.com
<br>강릉 박월준: woljun_park@kolon.com
<br>서산 김영민: youngmin_kim1@kolon.com
<br>SELEKT 김포 황도훈 : hangcome@kolon.com
'@
$ws.Range("D30").Value = $d30

# ---------------------------------------------------------------------------
# 3) D31 / D42: content unchanged — simply rewrite with the same value so
#    the stale <phoneticPr> attached to these shared strings is dropped
#    (matches the author's upload, which cleared phoneticPr everywhere).
# ---------------------------------------------------------------------------
$ws.Range("D31").Value = $ws.Range("D31").Value
$ws.Range("D42").Value = $ws.Range("D42").Value

# Row 42's height also shrank slightly in the re-save.
$ws.Rows.Item(42).RowHeight = 314.25

# ---------------------------------------------------------------------------
# 4) New row 43: "표시등 및 경고등" / HT203031
# ---------------------------------------------------------------------------
$d43 = @'
이 문서는 각 차량의 경고등이나 표시등에 대해 확인할 수 있습니다.
<br>
<br>표시등과 경고 시스템은 기능 켜짐, 시스템 작동 중, 결함 또는 심각한 오류 발생을 경고합니다.
<br>
<br>
<b3>적생 심벌</h3>
<br><img src="https://www.volvocars.com/images/support/img5605a7dcf7f527a8c0a8015253481e94_1_--_--_VOICEpnghigh.png" height="64" widgh="64"><string>경고</string><br>차량의 안전이나 주행성에 영향을 미칠 수 있는 결함이 감지되면 빨간색 경고등이 점등됩니다.<br>동시에 안내 메시지가 운전자 화면에 표시됩니다. 경고등은 다른 심볼과 함께 점등될 수도 있습니다.

'@

$ws.Range("A43").Value = "표시등 및 경고등"
$ws.Range("B43").Value = 45455
$ws.Range("C43").Value = "HT203031"
$ws.Range("D43").Value = $d43
$ws.Range("D43").WrapText = $true
$ws.Rows.Item(43).RowHeight = 174.75

# ---------------------------------------------------------------------------
# 5) Selection / active cell bookkeeping, matching the re-saved view state.
# ---------------------------------------------------------------------------
$ws.Range("C51").Select()

Write-Host "done"
